# Insert a new data row at row 116 (shifts existing rows 116-186 down to 117-187)
# and populate it with the new "Acelga" price record, per the commit:
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(116).Insert()

$ws.Cells.Item(116, 1).Value = 7
$ws.Cells.Item(116, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(116, 3).Value = "Ñuble"
$ws.Cells.Item(116, 4).Value = 44529
$ws.Cells.Item(116, 5).Value = 16
$ws.Cells.Item(116, 6).Value = 100112009
$ws.Cells.Item(116, 7).Value = "Acelga"
$ws.Cells.Item(116, 8).Value = "Sin especificar"
$ws.Cells.Item(116, 9).Value = "Primera"
$ws.Cells.Item(116, 10).Value = 100
$ws.Cells.Item(116, 11).Value = 350
$ws.Cells.Item(116, 12).Value = 400
$ws.Cells.Item(116, 13).Value = 375
$ws.Cells.Item(116, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(116, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(116, 16).Value = 375
$ws.Cells.Item(116, 17).Value = 1
$ws.Cells.Item(116, 18).Value = "Hortaliza"
